# "A lot of fixes"
#
# 1) Tighten up the page margins: right margin in from 851 twips (42.55pt)
#    down to 567 twips (28.35pt), and left margin out from 567 twips
#    (28.35pt) to 1134 twips (56.7pt). Top/bottom stay at 567 twips.
#    Word's PageSetup.*Margin properties are expressed in points, and
#    1 point == 20 twips, so we convert accordingly.
$d = $word.ActiveDocument

$d.PageSetup.TopMargin    = 567  / 20   # unchanged: 28.35pt / 567 twips
$d.PageSetup.RightMargin  = 567  / 20   # was 851 twips -> now 567 twips
$d.PageSetup.BottomMargin = 567  / 20   # unchanged: 28.35pt / 567 twips
$d.PageSetup.LeftMargin   = 1134 / 20   # was 567 twips  -> now 1134 twips

# 2) Wire up an APA bibliography source list for this document (this is
#    what Word's References > Citations & Bibliography > Manage Sources
#    writes out as customXml/item1.xml + its itemProps1.xml companion
#    part, with the style bookkeeping on the Bibliography object).
$d.Bibliography.BibliographyStyle = "APA"

$sourcesXml = '<?xml version="1.0" encoding="utf-8"?>' + "`r`n" + `
    '<b:Sources xmlns:b="http://schemas.openxmlformats.org/officeDocument/2006/bibliography" ' + `
    'xmlns="http://schemas.openxmlformats.org/officeDocument/2006/bibliography" ' + `
    'SelectedStyle="\APASixthEditionOfficeOnline.xsl" StyleName="APA" Version="6"/>'

$d.Bibliography.Sources.Add($sourcesXml) | Out-Null
$d.CustomXMLParts.Add($sourcesXml) | Out-Null
Write-Host "done"
